$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Record the Actual Result for TC03 (row 6)
$ws.Range("M6").Value = "Application does not run"

# Fill in the new Test Case (row 4 / TC01) details that were previously left blank:
#   G4 = Test Steps, I4 = Test Data, J4 = Test Coverage, K4 = Post Conditions
$ws.Range("G4").Value = "1.Open application" + [char]10 + "2.enter input values" + [char]10 + "3.start the prediction process"
$ws.Range("I4").Value = "Input values in expected range" + [char]10 + "[0,42,52812.09301,15609.38091,138961.2505]"
$ws.Range("J4").Value = "Statement"
$ws.Range("K4").Value = "Prediction shown on screen"

# Widen the "Test Data" column (I) so the newly entered text fits better
$ws.Columns.Item(9).ColumnWidth = 27.16666666666667

# Remove the trailing blank row that is no longer needed
$ws.Rows.Item(16).Delete()

# Leave the selection where the user last clicked after entering the data
$ws.Range("G5").Select()
